# Append two new game rows (16 and 17) to the "Spreads" sheet,
# matching the data produced by the latest scrape/grading run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Spreads")

$row16 = @{
    1 = 'Virginia Tech Hokies vs. Saint Joseph''s Hawks'
    2 = 'Nov 12 07:00PM ET'
    3 = 'Saint Joseph''s Hawks'
    4 = 3
    5 = 13
    6 = 8.25
    7 = 11
    8 = 0.5638787496655561
    9 = 13
    10 = 0.0400692258560322
    11 = 1.9
    12 = 7.5
    13 = 9
    14 = 9.9
    15 = 5.6
    16 = 0.1636200206280591
    17 = 612
    18 = 0.1342925608618203
    19 = 0.0232259319581084
    20 = 0.032
    21 = 0.24
    22 = 0.18
    23 = 0.192947480394298
    24 = 3
    25 = 156.5
    26 = 151.967978382876
    27 = 154.5
    28 = 156.5
    29 = 2.3
    30 = 151
    31 = 155
    32 = 152.935956765752
    33 = 149.68
    34 = 0.4699349568237442
    35 = 0.5300650431762558
    36 = -0.0538745669857796
    37 = 0.0062555193667319
    38 = 1
    39 = 1
    40 = 0
    41 = 1
}

$row17 = @{
    1 = 'Hawai''i Rainbow Warriors vs. Miss Valley St Delta Devils'
    2 = 'Nov 13 12:00AM ET'
    3 = 'Miss Valley St Delta Devils'
    4 = 3
    5 = 41
    6 = 32.85
    7 = 37.5
    8 = 0.6108053684643127
    9 = 41.25
    10 = 0.0813936037584303
    11 = 3.8
    12 = 31.7
    13 = 34
    14 = 39.1
    15 = 30.6
    18 = 0
    19 = 0
    20 = 0.005
    21 = 0.01
    22 = 0.0009999999999998
    23 = 0.0007367866889574
    24 = 3
    25 = 146.5
    26 = 144
    27 = 145.5
    28 = 146.5
    29 = 3
    30 = 144
    31 = 144
    32 = 148.374941716665
    33 = 140.98
    34 = 0.4849567874384211
    35 = 0.5150432125615789
    36 = -0.0388527363711027
    37 = -0.0064878879168899
    38 = 1
    39 = 1
    40 = 0
    41 = 0
}

foreach ($col in $row16.Keys) {
    $ws.Cells.Item(16, $col).Value = $row16[$col]
}

foreach ($col in $row17.Keys) {
    $ws.Cells.Item(17, $col).Value = $row17[$col]
}
